$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F1:G1").EntireColumn.Delete()
